# Append new order rows to the "Orders" sheet and refresh the encoded
# Number-summary string on the "Summary" sheet.
#
# Source data is exported with purely-numeric-looking cells (PackageID
# group markers in column A, and every "Number" value in column F) stored
# as TEXT rather than as numbers. Excel's COM layer auto-converts a plain
# numeric-looking string assigned to .Value into a real number, so for
# those cells we first flip NumberFormat to "@" (Text) to force the
# literal string to be kept as text instead of being parsed/rounded as a
# numeric value (this matters a lot for Summary!G2, which is a ~75 digit
# string that would otherwise be corrupted by float rounding).

$wb = $excel.ActiveWorkbook
$orders = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New rows to append after the existing data (rows 2-31), starting at row 32.
# Columns: A = PackageID (group marker, only set on the first row of a group),
#          C = FlowerName, F = Number.
$newRows = @(
    @{ A = "6"; C = "229_黄蝴蝶_Yellow Butterfly_Rosa rugosa Thunb._10stems"; F = "11" },
    @{ A = $null; C = "640_红辣椒_undefined_undefined_1bunch"; F = "5" },
    @{ A = $null; C = "512_松虫草粉_scabiosa pink_undefined_1bunch"; F = "10" },
    @{ A = $null; C = "419_松虫草红_scabiosa watermelon_undefined_1bunch"; F = "10" },
    @{ A = $null; C = "314_松虫草花边黑_scabiosa_undefined_1bunch"; F = "5" },
    @{ A = "1"; C = "3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "10" },
    @{ A = $null; C = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "10" },
    @{ A = $null; C = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"; F = "13" },
    @{ A = $null; C = "411_紫罗兰白_violet white_undefined_1bunch"; F = "10" },
    @{ A = $null; C = "327_文竹_asparagus fern_undefined_1bunch"; F = "10" },
    @{ A = $null; C = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"; F = "20" },
    @{ A = $null; C = "439_九星叶_undefined_undefined_1bunch"; F = "10" },
    @{ A = "2"; C = "320_雪柳花_Spiraea flower white_undefined_1bunch"; F = "5" },
    @{ A = $null; C = "586_洋牡丹白_undefined_undefined_1bunch"; F = "6" },
    @{ A = $null; C = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"; F = "30" },
    @{ A = $null; C = "113_绣球安娜绿_Hydrangea Anna Green_Hydrangea L._1stem"; F = "40" }
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $data = $newRows[$i]

    if ($null -ne $data.A) {
        $cellA = $orders.Cells.Item($rowNum, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $data.A
    }

    $orders.Cells.Item($rowNum, 3).Value = $data.C

    $cellF = $orders.Cells.Item($rowNum, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $data.F
}

# Rebuild the concatenated "Number" digit string (column F, rows 2..last)
# that is mirrored into Summary!G2. Use string interpolation (not "+") so
# PowerShell doesn't coerce the operands into numeric addition.
$lastRow = $startRow + $newRows.Count - 1
$numberString = "0"
for ($r = 2; $r -le $lastRow; $r++) {
    $cellValue = $orders.Cells.Item($r, 6).Value()
    $numberString = "$numberString$cellValue"
}

$g2 = $summary.Range("G2")
$g2.NumberFormat = "@"
$g2.Value = $numberString
